$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the old "800 metros" women's sheet entirely; keep the men's sheet and
# turn it into the new "4x100 metros" women's sheet (content gets overwritten below).
$wb.Worksheets.Item("Ciclos - F - 800 metros").Delete() | Out-Null

$ws = $wb.Worksheets.Item(1)
$ws.Name = "Ciclos - F - 4_100 metros"
$ws.Activate()

$ws.Range("A3").Value = "4*100 metros    FEMENINA"

$ws.Range("A7").Value = "Andrea, Martín"
$ws.Range("B7").Value = "1DAW"

$ws.Range("A8").Value = "Laura, Gómez"
$ws.Range("B8").Value = "1DAW"

$ws.Range("A9").Value = "María, López"
$ws.Range("B9").Value = "1DAW"

$ws.Range("A10").Value = "Paula, Díaz"
$ws.Range("B10").Value = "1DAW"
